$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 20.759853416306271
$ws.Range("C2").Value = 16.490162510701225
$ws.Range("D2").Value = 22.053182282657474
$ws.Range("E2").Value = 22.311463674313131

$ws.Range("B3").Value = 14.383597367489955
$ws.Range("C3").Value = 20.223775271097054
$ws.Range("D3").Value = 13.370704693699167
$ws.Range("E3").Value = 23.734676842306442

$ws.Range("B1:E3").Select()
